# Daily attendance processing - 2026-01-12 08:44:54
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Swap "System, <email>" -> "<email>, System" everywhere it occurs (Recorded By column) ---
$used = $ws.UsedRange
[void]$used.Replace("System, dnasr281@gmail.com", "dnasr281@gmail.com, System", 1)

# --- 2. Update summary statistics (Missing / Pending sessions) ---
$ws.Range("L7").Value = 33
$ws.Range("L8").Value = 54

# --- 3. Update per-group Missing/Pending counts for rows 15-20 ---
$ws.Range("P15").Value = 3
$ws.Range("Q15").Value = 3

$ws.Range("P16").Value = 2
$ws.Range("Q16").Value = 3

$ws.Range("P17").Value = 2
$ws.Range("Q17").Value = 3

$ws.Range("P18").Value = 2
$ws.Range("Q18").Value = 3

$ws.Range("P19").Value = 2
$ws.Range("Q19").Value = 3

$ws.Range("P20").Value = 3
$ws.Range("Q20").Value = 3

# --- 4. Rows that flipped from "Pending" (yellow) to "Not Recorded" (pink/red) ---
# Copy the formatting already used for "Not Recorded" rows (e.g. row 3) onto A:I
# of each affected row, then update the status text.
$fmtSource = $ws.Range("A3:I3")
$fmtSource.Copy()

$pendingRows = @(24, 50, 76, 102, 128, 154)
foreach ($r in $pendingRows) {
    $target = $ws.Range("A" + $r + ":I" + $r)
    $target.PasteSpecial(-4122)
    $ws.Range("I" + $r).Value = "Not Recorded"
}

$excel.CutCopyMode = 0
